$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All edited cells in columns D/E hold text-formatted values (e.g. "27.287.01",
# "  -2.34%  ") in the source sheet. Excel would otherwise auto-coerce these to
# numbers (losing trailing zeros, multi-dot grouping, etc.), so force the number
# format to Text ("@") before writing each value, matching the original cells.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.287.01"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -2.34%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.708.88"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -2.06%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "223.88"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.89%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5334"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -2.16%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2659"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -4.36%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06583"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.65%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.81"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -4.48%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07621"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -2.25%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.564"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -3.06%  "
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.713.56"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.81%  "
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.945.17"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.97%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5732"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -4.28%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅8169"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -3.00%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.70"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.91%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "27.294.03"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -2.33%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "215.76"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -4.00%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.664"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -3.98%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.45"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -4.44%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.968"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -4.57%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "141.83"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -3.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.740"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +4.21%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1213"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.86%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.268"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.86%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -5.45%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05395"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -4.94%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.292"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.88%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.493"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -5.83%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.424"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.93%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.641"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.60%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.49%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9487"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -3.32%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.408"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.52%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5857"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.20%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01627"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.71%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.860"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.49%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.043.61"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.43%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8411"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.18%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.79"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.36%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.851.16"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.03%  "
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "57.97"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -3.22%  "
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0₈109"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -6.76%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4506"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.55%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.003"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.09%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.086"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.73%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.81%  "
